$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Platform cell (C9) "Google Scholar" -> "Google" (new shared string)
$ws.Range("C9").Value = "Google"

# New Row 10 data
$ws.Range("A10").Value = 9

# B10 date — copy the date number format from B9 (style 1) then set the value
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 44202

$ws.Range("C10").Value = "Google"
$ws.Range("D10").Value = "MMR car price"

$ws.Range("E10").Value = "https://publish.manheim.com/en/help/mmr/mmr-under-the-hood.html"

# Add the hyperlink, then re-apply the existing "Link" cell style so the
# stored cell reuses the workbook's existing hyperlink style (matching E5).
$ws.Hyperlinks.Add($ws.Range("E10"), "https://publish.manheim.com/en/help/mmr/mmr-under-the-hood.html", "", "", "https://publish.manheim.com/en/help/mmr/mmr-under-the-hood.html")
$ws.Range("E10").Style = "Link"
